$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1. Prepare 6 extra pristine copies of the (still-empty) "Sheet1" sheet
#    so every new sheet starts from a clean slate (no custom cols, no
#    leftover selection) just like the original "Sheet1" tab.
# ----------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item(2)
$prev = $srcSheet
for ($i = 0; $i -lt 6; $i++) {
    $srcSheet.Copy($null, $prev)
    $prev = $wb.Worksheets.Item($prev.Index + 1)
}

# Rename the sheets into their final names, in order.
$names = @("ManageProductSearch","ManageCategoryNew","ManageProductNew","ManageFooter","ManageContact","ManageNewsNew","AdminUsers")
for ($i = 0; $i -lt $names.Count; $i++) {
    $wb.Worksheets.Item($i + 2).Name = $names[$i]
}

$wsSearch   = $wb.Worksheets.Item("ManageProductSearch")
$wsCategory = $wb.Worksheets.Item("ManageCategoryNew")
$wsProduct  = $wb.Worksheets.Item("ManageProductNew")
$wsFooter   = $wb.Worksheets.Item("ManageFooter")
$wsContact  = $wb.Worksheets.Item("ManageContact")
$wsNews     = $wb.Worksheets.Item("ManageNewsNew")
$wsAdmin    = $wb.Worksheets.Item("AdminUsers")

# ----------------------------------------------------------------------
# 2. ManageProductSearch
# ----------------------------------------------------------------------
$wsSearch.Cells.Item(1,1).Value = "Title"
$wsSearch.Cells.Item(1,2).Value = "ProductCode"
$wsSearch.Cells.Item(2,1).Value = "Curd12"
$wsSearch.Cells.Item(2,2).Value = "P1015"
$wsSearch.Cells.Item(3,1).Value = "Milk"
$wsSearch.Cells.Item(3,2).Value = "P992"
$wsSearch.Columns.Item(2).ColumnWidth = 13.33
$wsSearch.Range("B2").Select() | Out-Null

# ----------------------------------------------------------------------
# 3. ManageCategoryNew
# ----------------------------------------------------------------------
$wsCategory.Cells.Item(1,1).Value = "Category Value"
$wsCategory.Cells.Item(2,1).Value = "Discount"
$wsCategory.Columns.Item(1).ColumnWidth = 14.67
$wsCategory.Range("A2").Select() | Out-Null

# ----------------------------------------------------------------------
# 4. ManageProductNew
#    (shared-string insertion order observed in target: A1, B1, D1, C1)
# ----------------------------------------------------------------------
$wsProduct.Cells.Item(1,1).Value = "Title Value"
$wsProduct.Cells.Item(1,2).Value = "max qty"
$wsProduct.Cells.Item(1,4).Value = "stock"
$wsProduct.Cells.Item(1,3).Value = "price"
$wsProduct.Cells.Item(2,1).Value = "Milk"
$wsProduct.Cells.Item(2,2).Value = 45
$wsProduct.Cells.Item(2,3).Value = 150
$wsProduct.Cells.Item(2,4).Value = 45
$wsProduct.Columns.Item(1).ColumnWidth = 15.0
$wsProduct.Columns.Item(4).ColumnWidth = 12.0
$wsProduct.Range("E2").Select() | Out-Null

# ----------------------------------------------------------------------
# 5. ManageFooter
#    (shared-string insertion order observed in target: A1, A2, B2, B1, C1)
# ----------------------------------------------------------------------
$wsFooter.Cells.Item(1,1).Value = "Address"
$wsFooter.Cells.Item(2,1).Value = "Asiatic business center,Technopark Phase three,Trivandrum"
$wsFooter.Cells.Item(2,2).Value = "automationtesting@gmail.com"
$wsFooter.Cells.Item(1,2).Value = "Email ID"
$wsFooter.Cells.Item(1,3).Value = "Phone Number"
$wsFooter.Cells.Item(2,3).Value = 9876543210
$wsFooter.Range("A2:B2").WrapText = $true
$wsFooter.Range("A2:B2").HorizontalAlignment = -4131
$wsFooter.Rows.Item(2).RowHeight = 41.4
$wsFooter.Columns.Item(1).ColumnWidth = 21.83
$wsFooter.Columns.Item(2).ColumnWidth = 15.67
$wsFooter.Columns.Item(3).ColumnWidth = 13.0
$wsFooter.Range("A2").Select() | Out-Null

# ----------------------------------------------------------------------
# 6. ManageContact
#    (shared-string insertion order observed in target: A1(reuse), B1,
#     B2, C2, D1, E1, then C1(reuse) and the plain numbers)
# ----------------------------------------------------------------------
$wsContact.Cells.Item(1,1).Value = "Phone Number"
$wsContact.Cells.Item(1,2).Value = "eMail"
$wsContact.Cells.Item(2,2).Value = "def@gmail.com"
$wsContact.Cells.Item(2,3).Value = "abc houseDiv 14, Kadavanthra, Panampilly Nagar,Ernakulam"
$wsContact.Cells.Item(1,4).Value = "Delivery Time "
$wsContact.Cells.Item(1,5).Value = "Delivery Charge Limit"
$wsContact.Cells.Item(1,3).Value = "Address"
$wsContact.Cells.Item(2,1).Value = 2255
$wsContact.Cells.Item(2,4).Value = 60
$wsContact.Cells.Item(2,5).Value = 100
$wsContact.Range("C2").WrapText = $true
$wsContact.Rows.Item(2).RowHeight = 34.8
$wsContact.Columns.Item(1).ColumnWidth = 12.17
$wsContact.Columns.Item(2).ColumnWidth = 13.5
$wsContact.Columns.Item(3).ColumnWidth = 25.33
$wsContact.Columns.Item(4).ColumnWidth = 10.67
$wsContact.Columns.Item(5).ColumnWidth = 17.5
$wsContact.Range("A2").Select() | Out-Null

# ----------------------------------------------------------------------
# 7. ManageNewsNew
#    (shared-string insertion order observed in target: A2, A1)
# ----------------------------------------------------------------------
$wsNews.Cells.Item(2,1).Value = "Automation Today"
$wsNews.Cells.Item(1,1).Value = "News"
$wsNews.Columns.Item(1).ColumnWidth = 15.67
$wsNews.Range("L24").Select() | Out-Null

# ----------------------------------------------------------------------
# 8. AdminUsers
# ----------------------------------------------------------------------
$wsAdmin.Cells.Item(1,1).Value = "UserName"
$wsAdmin.Cells.Item(1,2).Value = "Password"
$wsAdmin.Cells.Item(2,1).Value = "Obsqura2"
$wsAdmin.Cells.Item(2,2).Value = 214356
$wsAdmin.Range("B2").Select() | Out-Null
$wsAdmin.Activate() | Out-Null

# ----------------------------------------------------------------------
# 9. LoginDetails selection tweak (loses tab-selected state, selection
#    moves to A2)
# ----------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("LoginDetails")
$wsLogin.Range("A2").Select() | Out-Null

# Re-activate AdminUsers last so it ends up the active/visible tab.
$wsAdmin.Activate() | Out-Null
$wsAdmin.Range("B2").Select() | Out-Null
